# Add "intermediate" prior blocks to the priors sheet, pushing the
# existing "strong" prior blocks down to make room (mirrors the
# "add deploy script from ignore" commit: new R_intermedMa_prior /
# R_intermedI0_prior blocks inserted above the existing
# R_strongMa_prior / R_strongI0_prior blocks, which now live at rows 14-20).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("priors")

# --- First, re-create the current "strong" blocks (rows 6-12) further
#     down the sheet, at rows 14-20, preserving their original content.

# R_strongMa_prior block -> rows 14-16
$ws.Range("A14").Value = "R_strongMa_prior"
$ws.Range("B14").Value = "Time-Series/Cumulative"
$ws.Range("C14").Value = "I0"
$ws.Range("D14").Value = "Uniform(1, 10)"

$ws.Range("C15").Value = "r1/ma2"
$ws.Range("D15").Value = "Beta(100, 900)"

$ws.Range("C16").Value = "ma2"
$ws.Range("D16").Value = "Beta(500, 500)"

# R_strongI0_prior block -> rows 18-20
$ws.Range("A18").Value = "R_strongI0_prior"
$ws.Range("B18").Value = "Time-Series/Cumulative"
$ws.Range("C18").Value = "I0"
$ws.Range("D18").Value = "Log-Normal(0.69, 0.05)"

$ws.Range("C19").Value = "r1"
$ws.Range("D19").Value = "Log-Normal(0, 5)"

$ws.Range("C20").Value = "ma2"
$ws.Range("D20").Value = "Uniform(0, 1)"

# --- Now overwrite rows 6-12 with the new "intermediate" prior blocks.
# The new string values are entered in this particular order so the
# shared-string table is populated the same way it was authored
# (names first, then the new numeric priors).

# Names of the two new blocks.
$ws.Range("A6").Value = "R_intermedMa_prior"
$ws.Range("A10").Value = "R_intermedI0_prior"

# R_intermedI0_prior new prior value.
$ws.Range("D10").Value = "Log-Normal(0.69, 0.5)"

# R_intermedMa_prior new prior values.
$ws.Range("D7").Value = "Beta(10, 90)"
$ws.Range("D8").Value = "Beta(50, 50)"

# Remaining (unchanged/reused) values for the two new blocks.
$ws.Range("B6").Value = "Time-Series/Cumulative"
$ws.Range("C6").Value = "I0"
$ws.Range("D6").Value = "Uniform(1, 10)"
$ws.Range("C7").Value = "r1/ma2"
$ws.Range("C8").Value = "ma2"

$ws.Range("B10").Value = "Time-Series/Cumulative"
$ws.Range("C10").Value = "I0"
$ws.Range("C11").Value = "r1"
$ws.Range("D11").Value = "Log-Normal(0, 5)"
$ws.Range("C12").Value = "ma2"
$ws.Range("D12").Value = "Uniform(0, 1)"

# Update the selected cell shown in the sheet view.
$ws.Range("D9").Select()
